$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2
$ws.Range("B2").Value = 9
$ws.Range("C2").Value = 17
$ws.Range("D2").Value = "sdvkjdsnvkj edited"
$ws.Range("E2").Value = "fsdklnvdsklvn"
$ws.Range("F2").Value = "knsdkjvndsk"
$ws.Range("G2").Value = "nvkjwnvkjsd"
$ws.Range("I2").Value = "knfdklvnsdl"

# Update row 3
$ws.Range("B3").Value = 10
$ws.Range("C3").Value = 17
$ws.Range("D3").Value = "sdvkjdsnvkjds"
$ws.Range("E3").Value = "nvsbkdjbvsnk"
$ws.Range("F3").Value = "sdvnkjvnsdk"
$ws.Range("G3").Value = "kjsdbvkdsj"
$ws.Range("I3").Value = "nvfksjdnvs"

# Delete row 4 entirely (shifts cells up / removes the row)
$ws.Rows.Item(4).Delete()
